# Applies the "new scenario analysis, changes to death equations" edit:
# adds prompt/delayed treatment-effect parameter rows for aph, ur (maternal
# death) and eclampsia, aph, ur (stillbirth), splitting them out from the
# existing combined parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# PASS 1: "_md" (maternal death) prompt/delayed rows
# ---------------------------------------------------------------------

# 1a. aph_treatment_effect_md lives at row 132 -> insert 2 rows after it
#     (i.e. at rows 133-134) for aph_prompt_treatment_effect_md /
#     aph_delayed_treatment_effect_md.
$ws.Rows.Item(133).Insert()
$ws.Rows.Item(133).Insert()

$ws.Range("A133").Value = "aph_prompt_treatment_effect_md"
$ws.Range("B133").Value = 0.5

$ws.Range("A134").Value = "aph_delayed_treatment_effect_md"
$ws.Range("B134").Value = 0.75

# 1b. ur_treatment_effect_md now lives at row 140 -> insert 2 rows after it
#     (i.e. at rows 141-142) for ur_prompt_treatment_effect_md /
#     ur_delayed_treatment_effect_md.
$ws.Rows.Item(141).Insert()
$ws.Rows.Item(141).Insert()

$ws.Range("A141").Value = "ur_prompt_treatment_effect_md"
$ws.Range("B141").Value = 0.5
$ws.Range("C141").Value = "DUMMY"

$ws.Range("A142").Value = "ur_delayed_treatment_effect_md"
$ws.Range("B142").Value = 0.75
$ws.Range("C142").Value = "DUMMY"

# ---------------------------------------------------------------------
# PASS 2: "_sb" (stillbirth) prompt/delayed rows
# ---------------------------------------------------------------------

# 2a. eclampsia_treatment_effect_sb now lives at row 131 -> insert 2 rows
#     before it (i.e. at rows 131-132) for eclampsia_prompt_treatment_effect_sb
#     / eclampsia_delayed_treatment_effect_sb.
$ws.Rows.Item(131).Insert()
$ws.Rows.Item(131).Insert()

$ws.Range("A131").Value = "eclampsia_prompt_treatment_effect_sb"
$ws.Range("B131").Value = 0.5

$ws.Range("A132").Value = "eclampsia_delayed_treatment_effect_sb"
$ws.Range("B132").Value = 0.75

# 2b. aph_treatment_effect_sb now lives at row 138 -> insert 2 rows after it
#     (i.e. at rows 139-140) for aph_prompt_treatment_effect_sb /
#     aph_delayed_treatment_effect_sb.
$ws.Rows.Item(139).Insert()
$ws.Rows.Item(139).Insert()

$ws.Range("A139").Value = "aph_prompt_treatment_effect_sb"
$ws.Range("B139").Value = 0.5

$ws.Range("A140").Value = "aph_delayed_treatment_effect_sb"
$ws.Range("B140").Value = 0.75

# 2c. ur_treatment_effect_sb now lives at row 147 -> insert 2 rows after it
#     (i.e. at rows 148-149) for ur_prompt_treatment_effect_sb /
#     ur_delayed_treatment_effect_sb.
$ws.Rows.Item(148).Insert()
$ws.Rows.Item(148).Insert()

$ws.Range("A148").Value = "ur_prompt_treatment_effect_sb"
$ws.Range("B148").Value = 0.5

$ws.Range("A149").Value = "ur_delayed_treatment_effect_sb"
$ws.Range("B149").Value = 0.75

# ---------------------------------------------------------------------
# Match the view state recorded in the workbook after the edit (selection
# moved to A131; the sheet was scrolled down to show the new rows).
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("A117"), $true)
$ws.Range("A131").Select()
